$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Bureaucrat" from C10 to D22, and "Gardens" from C16 to D21
$ws.Range("D21").Value2 = $ws.Range("C16").Value2
$ws.Range("D22").Value2 = $ws.Range("C10").Value2

$ws.Range("C10").ClearContents()
$ws.Range("C16").ClearContents()

# Update the active selection
$ws.Range("C19").Select()
